$wb = $excel.ActiveWorkbook
$st = $wb.Styles.Add("MyDateStyle")
$st.NumberFormat = "m/d/yyyy"
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("H3").Value = 42405
$ws1.Range("H3").Style = "MyDateStyle"
